$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.591.84"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "1.744.21"
$ws.Range("E3").Value = "  +1.09%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.51"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4922"
$ws.Range("E7").Value = "  +2.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2679"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06287"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("D10").Value = "1.743.08"
$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07044"
$ws.Range("E11").Value = "  -1.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.75"
$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6156"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.585"
$ws.Range("E14").Value = "  +0.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.14"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "26.592.17"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.57"
$ws.Range("E20").Value = "  -1.11%  "

$ws.Range("D21").Value = "1.966.18"
$ws.Range("E21").Value = "  +0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.578"
$ws.Range("E22").Value = "  +0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.722"
$ws.Range("E23").Value = "  -1.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.266"
$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.53"
$ws.Range("E25").Value = "  +2.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.47"
$ws.Range("E26").Value = "  +0.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.428"
$ws.Range("E27").Value = "  +1.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.764"
$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.77"
$ws.Range("E29").Value = "  +0.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.050"
$ws.Range("E30").Value = "  +2.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08037"
$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.744"
$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04630"
$ws.Range("E33").Value = "  +1.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9997"
$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.018"
$ws.Range("E36").Value = "  +3.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6384"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.071"
$ws.Range("E38").Value = "  +4.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8975"
$ws.Range("E39").Value = "  -4.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.424"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("E42").Value = "  +0.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.89"
$ws.Range("E43").Value = "  -5.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.431"
$ws.Range("E44").Value = "  -3.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3923"
$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.872"
$ws.Range("E46").Value = "  -1.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1183"
$ws.Range("E47").Value = "  -0.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05398"
$ws.Range("E48").Value = "  +1.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.55"
$ws.Range("E49").Value = "  -1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.775"
$ws.Range("E50").Value = "  -1.63%  "

$ws.Range("E51").Value = "  -0.11%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007306"
$ws.Range("E18").Value = "  +5.38%  "

$ws.Range("B19").Value = "BinanceUSD"
$ws.Range("C19").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.09%  "
